$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RSL(0)")
$ws.Range("A1").Value = "Command"
